# Update "datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 22:56"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6624227
$ws.Range("C4").Value = 34580
$ws.Range("D4").Value = 3891345
$ws.Range("E4").Value = 2535688
$ws.Range("G4").Value = 867
$ws.Range("H4").Value = 197194

# Row 11 - Sudafrica
$ws.Range("B11").Value = 646398
$ws.Range("C11").Value = 1960
$ws.Range("D11").Value = 574587
$ws.Range("E11").Value = 56433
$ws.Range("G11").Value = 113
$ws.Range("H11").Value = 15378

# Row 57 - Costa Rica
$ws.Range("B57").Value = 53969
$ws.Range("C57").Value = 1420
$ws.Range("D57").Value = 20710
$ws.Range("E57").Value = 32676
$ws.Range("G57").Value = 16
$ws.Range("H57").Value = 583

# Row 83 - Costa de Marfil
$ws.Range("B83").Value = 18916
$ws.Range("C83").Value = 47
$ws.Range("D83").Value = 17960
$ws.Range("E83").Value = 837

# Row 99 - Namibia
$ws.Range("B99").Value = 9437
$ws.Range("C99").Value = 181
$ws.Range("D99").Value = 4898
$ws.Range("E99").Value = 4441
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 98

# Rows 120-121 swap: Cuba and Cabo Verde exchange positions (Cabo Verde now
# ranks above Cuba), each keeping its own refreshed figures.
$ws.Range("A120").Value = "Cabo Verde"
$ws.Range("B120").Value = 4651
$ws.Range("C120").Value = 94
$ws.Range("D120").Value = 4076
$ws.Range("E120").Value = 531
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 44

$ws.Range("A121").Value = "Cuba"
$ws.Range("B121").Value = 4593
$ws.Range("C121").Value = 42
$ws.Range("D121").Value = 3844
$ws.Range("E121").Value = 643
$ws.Range("H121").Value = 106

# Row 122 - Ruanda
$ws.Range("B122").Value = 4534
$ws.Range("C122").Value = 55
$ws.Range("D122").Value = 2450
$ws.Range("E122").Value = 2062

# Rows 132-133 swap: Angola and Lituania exchange positions (Angola now
# ranks above Lituania), each keeping its own refreshed figures.
$ws.Range("A132").Value = "Angola"
$ws.Range("B132").Value = 3279
$ws.Range("C132").Value = 62
$ws.Range("D132").Value = 1288
$ws.Range("E132").Value = 1860
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 131

$ws.Range("A133").Value = "Lituania"
$ws.Range("B133").Value = 3243
$ws.Range("C133").Value = 44
$ws.Range("D133").Value = 2049
$ws.Range("E133").Value = 1108
$ws.Range("H133").Value = 86

# Row 142 - Sudan del Sur
$ws.Range("B142").Value = 2568
$ws.Range("C142").Value = 13
$ws.Range("E142").Value = 1229

# Row 156 - Togo
$ws.Range("B156").Value = 1548
$ws.Range("C156").Value = 11
$ws.Range("D156").Value = 1166
$ws.Range("E156").Value = 345

# Row 158 - Burkina Faso
$ws.Range("B158").Value = 1499
$ws.Range("C158").Value = 13
$ws.Range("D158").Value = 1127
$ws.Range("E158").Value = 316

# Row 168 - Santo Tome y Principe
$ws.Range("B168").Value = 906
$ws.Range("C168").Value = 5
$ws.Range("E168").Value = 25
